# Penalty Reward System update (unfinished) - update forecast values and
# shift week-start dates forward by one week on the "Forecast Comparison"
# sheet, and refresh the dependent metrics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Force columns/cells to be treated as plain text so the date-like strings
# and numeric-looking strings are not auto-converted to dates/numbers.
$ws1.Range("B2:B17").NumberFormat = "@"
$ws2.Range("B2:B15").NumberFormat = "@"

# --- Forecast Comparison sheet: Week_Start_Date (B) and MyForecast (D) ---

$ws1.Range("B2").Value = "2025-01-12"
$ws1.Range("D2").Value = 33

$ws1.Range("B3").Value = "2025-01-19"
$ws1.Range("D3").Value = 32

$ws1.Range("B4").Value = "2025-01-26"
$ws1.Range("D4").Value = 32

$ws1.Range("B5").Value = "2025-02-02"
$ws1.Range("D5").Value = 33

$ws1.Range("B6").Value = "2025-02-09"
$ws1.Range("D6").Value = 34

$ws1.Range("B7").Value = "2025-02-16"
$ws1.Range("D7").Value = 35

$ws1.Range("B8").Value = "2025-02-23"
$ws1.Range("D8").Value = 36

$ws1.Range("B9").Value = "2025-03-02"
$ws1.Range("D9").Value = 36

$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("D10").Value = 36

$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("D11").Value = 37

$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("D12").Value = 38

$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("D13").Value = 39

$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("D14").Value = 32

$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("D15").Value = 38

$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("D16").Value = 37

$ws1.Range("B17").Value = "2025-04-27"
$ws1.Range("D17").Value = 36

# --- Summary sheet: recomputed metrics ---

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-05"
$ws2.Range("B4").Value = "184"
$ws2.Range("B5").Value = "73"
$ws2.Range("B6").Value = "71"
$ws2.Range("B7").Value = "45"
$ws2.Range("B8").Value = "7854 units"
$ws2.Range("B9").Value = "566"
$ws2.Range("B10").Value = "270"
$ws2.Range("B11").Value = "129"
$ws2.Range("B12").Value = "39"
$ws2.Range("B13").Value = "2025-03-30"
$ws2.Range("B14").Value = "32"
$ws2.Range("B15").Value = "2025-01-19"
